$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (was row 15 before the edit)
$ws.Range("A10").Value = 111541121
$ws.Range("B10").Value = 79444
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1049
$ws.Range("F10").Value = "Kortskaftad ärgspik"
$ws.Range("G10").Value = "Microcalicium ahlneri"
$ws.Range("H10").Value = "Tibell"
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 693460.9606228607
$ws.Range("R10").Value = 6551521.405726598
$ws.Range("AF10").Value = ""
$ws.Range("AO10").Value = "silverstubbe av tall"

# Row 11 (was row 10 before the edit)
$ws.Range("A11").Value = 111541122
$ws.Range("B11").Value = 5112
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 102204
$ws.Range("F11").Value = "Grönhjon"
$ws.Range("G11").Value = "Callidium aeneum"
$ws.Range("H11").Value = "(De Geer, 1775)"
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = "äldre gnagspår"
$ws.Range("N11").Value = ""
$ws.Range("Q11").Value = 693344.0451535647
$ws.Range("R11").Value = 6551526.82974836
$ws.Range("AF11").Value = ""
$ws.Range("AO11").Value = "låga av gran"

# Row 12 (was row 16 before the edit)
$ws.Range("A12").Value = 111541129
$ws.Range("B12").Value = 5113
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 100526
$ws.Range("F12").Value = "Bronshjon"
$ws.Range("G12").Value = "Callidium coriaceum"
$ws.Range("H12").Value = "Paykull, 1800"
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "äldre gnagspår"
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 693328.6441019299
$ws.Range("R12").Value = 6551545.628735202
$ws.Range("AF12").Value = ""
$ws.Range("AO12").Value = "torrgran"

# Row 13 (was row 12 before the edit)
$ws.Range("A13").Value = 111541115
$ws.Range("B13").Value = 89405
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("Q13").Value = 693612.9253791923
$ws.Range("R13").Value = 6551435.326171798
$ws.Range("AF13").Value = ""
$ws.Range("AO13").Value = "låga av gran"

# Row 14 (was row 13 before the edit)
$ws.Range("A14").Value = 111541120
$ws.Range("B14").Value = 79444
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1049
$ws.Range("F14").Value = "Kortskaftad ärgspik"
$ws.Range("G14").Value = "Microcalicium ahlneri"
$ws.Range("H14").Value = "Tibell"
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("Q14").Value = 693513.2669972532
$ws.Range("R14").Value = 6551517.868690074
$ws.Range("AF14").Value = ""
$ws.Range("AO14").Value = "silverstubbe av tall"

# Row 15 (was row 11 before the edit)
$ws.Range("A15").Value = 111541128
$ws.Range("B15").Value = 5113
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 100526
$ws.Range("F15").Value = "Bronshjon"
$ws.Range("G15").Value = "Callidium coriaceum"
$ws.Range("H15").Value = "Paykull, 1800"
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = "färska gnagspår"
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 693570.8046739453
$ws.Range("R15").Value = 6551451.742365629
$ws.Range("AF15").Value = ""
$ws.Range("AO15").Value = "torrgran"

# Row 16 (was row 14 before the edit)
$ws.Range("A16").Value = 111541118
$ws.Range("B16").Value = 94851
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 2569
$ws.Range("F16").Value = "Stor revmossa"
$ws.Range("G16").Value = "Bazzania trilobata"
$ws.Range("H16").Value = "(L.) Gray"
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("Q16").Value = 693461.6376634488
$ws.Range("R16").Value = 6551559.049034445
$ws.Range("AF16").Value = ""
$ws.Range("AO16").Value = ""

